$d = $word.ActiveDocument
$d.Content.Find.Execute("Guassian", $false, $false, $false, $false, $false, $true, 1, $false, "Gaussian", 2)
